$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 302; this shifts existing rows 302-386 down to 303-387
$ws.Rows.Item(302).Insert()

# Populate the newly inserted row 302 with the new data entry
$ws.Range("A302").Value = 5
$ws.Range("B302").Value = "Macroferia Regional de Talca"
$ws.Range("C302").Value = "Maule"
$ws.Range("D302").Value = 44985
$ws.Range("E302").Value = 7
$ws.Range("F302").Value = 100112045
$ws.Range("G302").Value = "Zapallo"
$ws.Range("H302").Value = "Camote"
$ws.Range("I302").Value = "1a (cosecha)"
$ws.Range("J302").Value = 900
$ws.Range("K302").Value = 300
$ws.Range("L302").Value = 300
$ws.Range("M302").Value = 300
$ws.Range("N302").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O302").Value = "Región del Maule"
$ws.Range("P302").Value = 300
$ws.Range("Q302").Value = 1
$ws.Range("R302").Value = "Hortaliza"
